$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
  'dringen',
  'rasen',
  'löschen',
  'bauen',
  'posten',
  'saufen',
  'flehen',
  'süßen',
  'liegen',
  'mögen',
  'klettern',
  'pressen',
  'drehen',
  'führen',
  'streichen',
  'trauen',
  'tollen',
  'öffnen',
  'zeigen',
  'fangen',
  'garen',
  'decken',
  'schleppen',
  'schlafen',
  'locken',
  'backen',
  'enden',
  'warnen',
  'starten',
  'trennen',
  'quälen',
  'bellen',
  'beten',
  'schützen',
  'loben',
  'testen',
  'schultern',
  'ändern',
  'wohnen',
  'wechseln',
  'folgen',
  'stehlen',
  'wandern',
  'spannen',
  'messen',
  'feiern',
  'wehtun',
  'geben',
  'wüten',
  'hören',
  'husten',
  'arten',
  'schmecken',
  'meinen',
  'schrecken',
  'wehen',
  'kichern',
  'siegen',
  'münzen',
  'kümmern',
  'ehren',
  'fassen',
  'scheinen',
  'sterben',
  'krachen',
  'rauben',
  'hupen',
  'pflanzen',
  'sprengen',
  'helfen',
  'lockern',
  'trösten',
  'scheitern',
  'trotzen',
  'fallen',
  'kehren',
  'treten',
  'bitten',
  'schauen',
  'schalten',
  'kranken',
  'biegen',
  'hauen',
  'erben',
  'ärgern',
  'sorgen',
  'graben',
  'merken',
  'gründen',
  'schwingen',
  'zielen',
  'klagen',
  'stechen',
  'klingen',
  'jubeln',
  'fahren',
  'zögern',
  'regnen',
  'freuen',
  'rücken',
  'filmen',
  'formen',
  'greifen',
  'schlucken',
  'treiben',
  'schwächen',
  'lächeln',
  'teilen',
  'räumen',
  'lügen',
  'wundern',
  'heißen',
  'flüstern',
  'sichern',
  'kosten',
  'heulen',
  'rufen',
  'wachsen',
  'ziehen',
  'schließen',
  'hindern',
  'stecken',
  'bluten',
  'malen',
  'jagen',
  'seufzen',
  'spenden',
  'wirken',
  'albern',
  'achten',
  'hassen',
  'passen',
  'träumen',
  'brauchen',
  'zünden',
  'retten',
  'wenden',
  'runden',
  'stärken',
  'hoffen',
  'wecken',
  'spinnen',
  'platzen',
  'leugnen',
  'betteln',
  'werfen',
  'weichen',
  'bergen',
  'proben',
  'fällen',
  'töten',
  'dienen',
  'planen',
  'kennen',
  'tragen',
  'schulden',
  'heben',
  'brauen',
  'tropfen',
  'liefern',
  'irren',
  'schwören',
  'knarren',
  'spielen',
  'segeln',
  'äußern',
  'werden',
  'schenken',
  'sinken',
  'leisten',
  'gelten',
  'schreiten',
  'machen',
  'stammen',
  'mauern',
  'lesen',
  'haben',
  'stoppen',
  'spüren',
  'pfeifen',
  'altern',
  'heilen',
  'sparen',
  'fließen',
  'suchen',
  'sperren',
  'grüßen',
  'reizen',
  'kürzen',
  'flüchten',
  'boxen',
  'fischen'
)

for ($i = 0; $i -lt $words.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 2).Value = $words[$i]
}
